$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 26, pushing existing rows 26-68 down to 28-70.
$ws.Rows.Item(26).Resize(2).Insert()

# New row 26: Camote record dated 2022-04-11 (serial 44662), Perú origin.
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = "2022-04-11"
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 100114002
$ws.Cells.Item(26, 7).Value = "Camote"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 20
$ws.Cells.Item(26, 11).Value = 18000
$ws.Cells.Item(26, 12).Value = 18000
$ws.Cells.Item(26, 13).Value = 18000
$ws.Cells.Item(26, 14).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(26, 15).Value = "Perú"
$ws.Cells.Item(26, 16).Value = 1200
$ws.Cells.Item(26, 17).Value = 15
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# New row 27: Camote record dated 2022-04-11 (serial 44662), Perú origin.
$ws.Cells.Item(27, 1).Value = 10
$ws.Cells.Item(27, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(27, 3).Value = "La Araucanía"
$ws.Cells.Item(27, 4).Value = "2022-04-11"
$ws.Cells.Item(27, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(27, 5).Value = 9
$ws.Cells.Item(27, 6).Value = 100114002
$ws.Cells.Item(27, 7).Value = "Camote"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 50
$ws.Cells.Item(27, 11).Value = 18000
$ws.Cells.Item(27, 12).Value = 18000
$ws.Cells.Item(27, 13).Value = 18000
$ws.Cells.Item(27, 14).Value = "`$/malla 20 kilos"
$ws.Cells.Item(27, 15).Value = "Perú"
$ws.Cells.Item(27, 16).Value = 900
$ws.Cells.Item(27, 17).Value = 20
$ws.Cells.Item(27, 18).Value = "Hortaliza"
